$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = '<and the aptain before he lodged of one him him hontorcdely lodged of on him hundrition of one him him huntoly lodged.>'
$ws.Range("C2").Value = 59.72740126689256
$ws.Range("B3").Value = '<prior his regulard in accared in in the bingert in accared in accared in in three modely to careford in accared in in the breald separationd acciate accialed at ack.>'
$ws.Range("C3").Value = 59.81399736193821
$ws.Range("B4").Value = '<in then then two pored soupwards, in two pore ded two porace>'
$ws.Range("C4").Value = 61.15430567204474
$ws.Range("B5").Value = '<benevey stoped in the transle the gonman ran to the gonman ran waited in the transil and to the gunt and to the gor.>'
$ws.Range("C5").Value = 59.89389088374904
$ws.Range("B6").Value = '<and to as frequent to the business of the country on the busiant to assiation.>'
$ws.Range("C6").Value = 58.82777656093278
$ws.Range("B7").Value = '<he diring to get had every on apartain that he had every on apartal fraid every on apartaing to get curtail to gett.>'
$ws.Range("C7").Value = 59.40708229374385
$ws.Range("B8").Value = '<while approrat the hells the appror.>'
$ws.Range("C8").Value = 56.79272399326966
$ws.Range("B9").Value = '<toward any of the inights.>'
$ws.Range("C9").Value = 62.57341677915874
$ws.Range("B10").Value = '<in the morey in the marinest as to the conviction as to the conviction assist assist assistran.>'
$ws.Range("C10").Value = 58.34802043946114
$ws.Range("B11").Value = '<to the government, to the bitttorement, all and every the government, to the bittered the government>'
$ws.Range("C11").Value = 59.9437863283374
$ws.Range("B12").Value = '<but the manney manned accrea to they mannes, they neaver the less of making man accred to their perper been accred to their perper depers.>'
$ws.Range("C12").Value = 61.96117091740761
$ws.Range("B13").Value = '<in the case of convicted merderers of convicted merderers of convicted.>'
$ws.Range("C13").Value = 63.10442032634959
$ws.Range("B14").Value = '<thrance, catther, the cating and for the condemned at than for the colonay are a moval to deating execution.>'
$ws.Range("C14").Value = 60.38267898876224
$ws.Range("B15").Value = '<a greaking to in oswald states thing to the union oswald states thing to the union oswald states.>'
$ws.Range("C15").Value = 60.40272963374597
$ws.Range("B16").Value = '<he was himself act dearch, and all chere of the locaution an all chirced works.>'
$ws.Range("C16").Value = 60.94305764667524
$ws.Range("B17").Value = '<may not anot and sixty for sixty for.>'
$ws.Range("C17").Value = 60.19458891307603
$ws.Range("B18").Value = '<hose starce starce starce or exclosed loce by the staril mit>'
$ws.Range("C18").Value = 62.56909840138632
$ws.Range("B19").Value = '<he went one crediation at recause on cours on cocurs on cocurs at once on course on course on course a para credimented. and brokers a pose on crediation ate.>'
$ws.Range("C19").Value = 62.94524609804539
$ws.Range("B20").Value = '<and if the counsed if the counsed if the counsed if the treasury, of the should be mained for the may,>'
$ws.Range("C20").Value = 61.28472048555136
$ws.Range("B21").Value = '<some of the remainent in their place founds of the walls of the walls of the walls of the remain in their place founds of the walls of the walls.>'
$ws.Range("C21").Value = 61.12600614736264
$ws.Range("B22").Value = '<oswald was five five five five five five five five five nincasp slend.>'
$ws.Range("C22").Value = 62.18692870230829
$ws.Range("B23").Value = '<the service preferres preferres preferres preferres preferres.>'
$ws.Range("C23").Value = 63.47775454417791
$ws.Range("B24").Value = '<our lang with presisty hawels of a be down all sistem, however, quared now sisty the reprison.>'
$ws.Range("C24").Value = 60.09847064669823
$ws.Range("B25").Value = '<the game with was nearly a paine however once was nearly the gain how ever once may photogras walk, nearly field.>'
$ws.Range("C25").Value = 60.53341840069663
$ws.Range("B26").Value = '<when he handsmist with the was distorigst reate a was discrewards, rester.>'
$ws.Range("C26").Value = 66.36076053615014
$ws.Range("B27").Value = '<as was accaing of a dated made made mained sixty thress card, nineteen sixty,>'
$ws.Range("C27").Value = 57.50682394581035
$ws.Range("B28").Value = '<five point six precoide saffied in piort point saffied.>'
$ws.Range("C28").Value = 62.0867204129167
$ws.Range("B29").Value = '<oswald#s roove our overlo.>'
$ws.Range("C29").Value = 62.69042539122574
$ws.Range("B30").Value = '<loserlations warned him a befrind him a befrind him a befrined.>'
$ws.Range("C30").Value = 62.3216099750752
$ws.Range("B31").Value = '<this the metropomly appomates and whole appon the commors, and who strong in force and who strongate, appon the commonly appomect a pomect,>'
$ws.Range("C31").Value = 61.47040861236476
$ws.Range("B32").Value = '<and regarding the commission amotorig by the commission amotoring the commission regarding>'
$ws.Range("C32").Value = 57.85948019895099
$ws.Range("B33").Value = '<to her hus been the new orleans member husbands member husard chapter.>'
$ws.Range("C33").Value = 57.8845343960218
$ws.Range("B34").Value = '<who dispeated with the president#s the president#s then travelon, who divelm,>'
$ws.Range("C34").Value = 62.40244961798211
$ws.Range("B35").Value = '<his appiel from in compience appielf one the said broughts on the said brough the said brough hopse on the said brough hopse on the saffe up hopsent.>'
$ws.Range("C35").Value = 60.83797757135354
$ws.Range("B36").Value = '<addhost the bureau has unned such such such such such such such succe>'
$ws.Range("C36").Value = 60.56046421631492
$ws.Range("B37").Value = '<in the seventy the serippated with the seripped transpapers him the service papers him the served transpapers horw arrow hored with trancemed end quote.>'
$ws.Range("C37").Value = 59.71174780708824
$ws.Range("B38").Value = '<clied to supplace and rety palace arests of the forgoons of the four going res.>'
$ws.Range("C38").Value = 61.77054535291613
$ws.Range("B39").Value = '<to have a thousand pounds and frogreen previded thousand pounds and fraud tained preatence.>'
$ws.Range("C39").Value = 59.13387095955949
$ws.Range("B40").Value = '<and to crivids similary to crives.>'
$ws.Range("C40").Value = 61.74911431503838
$ws.Range("B41").Value = '<which handscroses was one speak hind quote.>'
$ws.Range("C41").Value = 60.95474569714644
$ws.Range("B42").Value = '<they talk a days were home of a more in are in are in are in are in are in are in are.>'
$ws.Range("C42").Value = 59.35193382285891
